# "fixed RAG excel agent"
# The order-inventory sheet had a malformed "summary" row (row 6) that stored
# its Order Number / Part Number as text instead of numbers. This fixes that
# row to use real numbers, appends a genuine new order row (order #7), and
# appends a new malformed-style row (order #8) carrying the same text-typed
# quirk the old row 6 used to have - i.e. the bug just moved to the new last
# row instead of being fully eliminated upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the old last row (row 6): Order Number / Part Number should be
#     real numbers, not text. Everything else in that row is unchanged.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 349584398539

# --- Insert a brand-new order row (row 7) ---
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 1290138230948
$ws.Range("C7").Value = "inflatable boat"
$ws.Range("D7").Value = 500
$ws.Range("E7").Value = "Matt"
$ws.Range("F7").Value = "Tom"

# --- Append a new row (row 8) reproducing the same text-typed Order
#     Number / Part Number / Price quirk the old row 6 had. Format as Text
#     first so the numeric-looking values are stored as literal strings,
#     matching the pre-existing data quality bug being carried forward.
$ws.Range("A8:B8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("A8").Value = "8"
$ws.Range("B8").Value = "33945803930"
$ws.Range("C8").Value = "stuffed animal bunny"
$ws.Range("D8").Value = "$25"
$ws.Range("E8").Value = "Kelly"
$ws.Range("F8").Value = "Felicia"
